$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '44.794.90'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  +1.43%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.246.51'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  +0.06%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.68'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  -0.15%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.55'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = '  -0.51%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.572'
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = '  -0.52%  '

$ws.Range("E8").Value = '  +0.22%  '

$ws.Range("E9").Value = '  -1.28%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.13'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  +0.53%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0805'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = '  -1.53%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.21'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = '  -0.73%  '

$ws.Range("E13").Value = '  -0.09%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.591.14'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  +0.12%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.329.54'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '  +3.77%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.840'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  +0.75%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.61'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  -0.24%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '44.601.54'
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0949'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = '  -2.55%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.97'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '  -2.31%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.27'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  -1.78%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.44'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = '  -0.03%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '239.59'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = '  +1.19%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.97'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = '  +0.35%  '

$ws.Range("E25").Value = '  -0.99%  '

$ws.Range("E26").Value = '  -0.10%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.30'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  +4.78%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.86'
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = '  -0.89%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '37.83'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = '  -0.69%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.04'
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = '  +0.21%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '19.97'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = '  -0.55%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '150.72'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = '  -0.88%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0796'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = '  -1.20%  '

$ws.Range("E34").Value = '  +1.06%  '

$ws.Range("E35").Value = '  -8.60%  '

$ws.Range("E36").Value = '  -1.10%  '

$ws.Range("E37").Value = '  -0.51%  '

$ws.Range("E38").Value = '  +4.41%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.13'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = '  +2.17%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.39'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = '  -0.33%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.78'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = '  -2.34%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0301'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = '  +1.04%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.848.86'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  +7.24%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.75'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  +13.89%  '

$ws.Range("B46").Value = 'BitcoinSV'

$ws.Range("C46").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '79.64'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = '  -4.69%  '

$ws.Range("B47").Value = 'Algorand'

$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.190'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = '  -0.52%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '98.86'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = '  -1.40%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.90'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = '  +0.98%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '69.07'
$ws.Range("D50").Style = "Normal"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '54.41'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = '  -0.49%  '
